$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 142860720
$ws.Cells.Item(64, 9).Value = 1000000000
$ws.Cells.Item(64, 10).Value = 4166.6665
$ws.Cells.Item(64, 11).Value = 1000000000
$ws.Cells.Item(64, 12).Value = 4166.6665
$ws.Cells.Item(64, 13).Value = -999999752
$ws.Cells.Item(64, 14).Value = -4662.6665

$ws.Cells.Item(67, 8).Value = 142860720
$ws.Cells.Item(67, 9).Value = 1000000000
$ws.Cells.Item(67, 10).Value = 4166.6665
$ws.Cells.Item(67, 11).Value = 1000000000
$ws.Cells.Item(67, 12).Value = 4166.6665
$ws.Cells.Item(67, 13).Value = -999999142
$ws.Cells.Item(67, 14).Value = -5882.6665

$ws.Cells.Item(74, 8).Value = 3214.4243
$ws.Cells.Item(74, 9).Value = 2454.889
$ws.Cells.Item(74, 10).Value = 3499.25
$ws.Cells.Item(74, 11).Value = 2454.889
$ws.Cells.Item(74, 12).Value = 3499.25
$ws.Cells.Item(74, 13).Value = -1518.889
$ws.Cells.Item(74, 14).Value = -5371.25

$ws.Cells.Item(77, 8).Value = 3214.4243
$ws.Cells.Item(77, 9).Value = 2454.889
$ws.Cells.Item(77, 10).Value = 3499.25
$ws.Cells.Item(77, 11).Value = 12274.445
$ws.Cells.Item(77, 12).Value = 17496.25
$ws.Cells.Item(77, 13).Value = -7594.445
$ws.Cells.Item(77, 14).Value = -26856.25

$ws.Cells.Item(100, 8).Value = 2302.8572
$ws.Cells.Item(100, 9).Value = 2064
$ws.Cells.Item(100, 10).Value = 2900
$ws.Cells.Item(100, 11).Value = 2064
$ws.Cells.Item(100, 12).Value = 2900
$ws.Cells.Item(100, 13).Value = -1523
$ws.Cells.Item(100, 14).Value = -3982

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 608.2222
$ws.Cells.Item(4, 9).Value = 659.375
$ws.Cells.Item(4, 10).Value = 199
$ws.Cells.Item(4, 11).Value = 659.375
$ws.Cells.Item(4, 12).Value = 199
$ws.Cells.Item(4, 13).Value = -543.375
$ws.Cells.Item(4, 14).Value = -431

$ws.Cells.Item(16, 8).Value = 1006
$ws.Cells.Item(16, 9).Value = 1006
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1006
$ws.Cells.Item(16, 12).Value = ""
$ws.Cells.Item(16, 13).Value = -719
$ws.Cells.Item(16, 14).Value = 0

$ws.Cells.Item(21, 8).Value = 1896
$ws.Cells.Item(21, 9).Value = 1896
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 1896
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = ""
$ws.Cells.Item(21, 14).Value = -1522

$ws.Cells.Item(58, 8).Value = 30000
$ws.Cells.Item(58, 10).Value = 30000
$ws.Cells.Item(58, 12).Value = 30000
$ws.Cells.Item(58, 14).Value = -30860

$ws.Cells.Item(97, 8).Value = 50757.95
$ws.Cells.Item(97, 9).Value = 59373.47
$ws.Cells.Item(97, 10).Value = 1936.6666
$ws.Cells.Item(97, 11).Value = 59373.47
$ws.Cells.Item(97, 12).Value = 1936.6666
$ws.Cells.Item(97, 13).Value = -58877.47
$ws.Cells.Item(97, 14).Value = -2928.6666

$ws.Cells.Item(102, 8).Value = 1457.6923
$ws.Cells.Item(102, 9).Value = 1217.7778
$ws.Cells.Item(102, 10).Value = 1997.5
$ws.Cells.Item(102, 11).Value = 1217.7778
$ws.Cells.Item(102, 12).Value = 1997.5
$ws.Cells.Item(102, 13).Value = 404.2221999999999
$ws.Cells.Item(102, 14).Value = -5241.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 29666.666
$ws.Cells.Item(82, 9).Value = 29000
$ws.Cells.Item(82, 10).Value = 30000
$ws.Cells.Item(82, 11).Value = 29000
$ws.Cells.Item(82, 12).Value = 30000
$ws.Cells.Item(82, 13).Value = -28617
$ws.Cells.Item(82, 14).Value = -30766

$ws.Cells.Item(85, 8).Value = 29666.666
$ws.Cells.Item(85, 9).Value = 29000
$ws.Cells.Item(85, 10).Value = 30000
$ws.Cells.Item(85, 11).Value = 29000
$ws.Cells.Item(85, 12).Value = 30000
$ws.Cells.Item(85, 13).Value = -27674
$ws.Cells.Item(85, 14).Value = -32652

$ws.Cells.Item(86, 8).Value = 2224.375
$ws.Cells.Item(86, 9).Value = 1452.7646
$ws.Cells.Item(86, 10).Value = 2794.6956
$ws.Cells.Item(86, 11).Value = 1452.7646
$ws.Cells.Item(86, 12).Value = 2794.6956
$ws.Cells.Item(86, 13).Value = -329.7646
$ws.Cells.Item(86, 14).Value = -5040.6956

$ws.Cells.Item(89, 8).Value = 2224.375
$ws.Cells.Item(89, 9).Value = 1452.7646
$ws.Cells.Item(89, 10).Value = 2794.6956
$ws.Cells.Item(89, 11).Value = 7263.823
$ws.Cells.Item(89, 12).Value = 13973.478
$ws.Cells.Item(89, 13).Value = -1647.823
$ws.Cells.Item(89, 14).Value = -25205.478

$ws.Cells.Item(94, 8).Value = 1050.9429
$ws.Cells.Item(94, 9).Value = 687.6667
$ws.Cells.Item(94, 10).Value = 1435.5883
$ws.Cells.Item(94, 11).Value = 687.6667
$ws.Cells.Item(94, 12).Value = 1435.5883
$ws.Cells.Item(94, 13).Value = -236.6667
$ws.Cells.Item(94, 14).Value = -2337.5883

$ws.Cells.Item(97, 8).Value = 18647.334
$ws.Cells.Item(97, 9).Value = 15000
$ws.Cells.Item(97, 11).Value = 15000
$ws.Cells.Item(97, 13).Value = -14009

$ws.Cells.Item(99, 8).Value = 998.3333
$ws.Cells.Item(99, 9).Value = 758.3333
$ws.Cells.Item(99, 10).Value = 1238.3334
$ws.Cells.Item(99, 11).Value = 758.3333
$ws.Cells.Item(99, 12).Value = 1238.3334
$ws.Cells.Item(99, 13).Value = 739.6667
$ws.Cells.Item(99, 14).Value = -4234.3334

$ws.Cells.Item(100, 8).Value = 48000
$ws.Cells.Item(100, 10).Value = 48000
$ws.Cells.Item(100, 12).Value = 48000
$ws.Cells.Item(100, 14).Value = -50164

$ws.Cells.Item(103, 8).Value = 22990.334
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 22990.334
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 12).Value = ""
$ws.Cells.Item(103, 13).Value = 22990.334
$ws.Cells.Item(103, 14).Value = -25334.334

$ws.Cells.Item(105, 8).Value = 2317.3
$ws.Cells.Item(105, 9).Value = 1410
$ws.Cells.Item(105, 10).Value = 2706.1428
$ws.Cells.Item(105, 11).Value = 1410
$ws.Cells.Item(105, 12).Value = 2706.1428
$ws.Cells.Item(105, 13).Value = 337
$ws.Cells.Item(105, 14).Value = -6200.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(100, 8).Value = 24390
$ws.Cells.Item(100, 10).Value = 24390
$ws.Cells.Item(100, 12).Value = 24390
$ws.Cells.Item(100, 14).Value = -26554

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 892.3542
$ws.Cells.Item(131, 10).Value = 1003.7561
$ws.Cells.Item(131, 12).Value = 3011.2683
$ws.Cells.Item(131, 14).Value = -13091.2683

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 1000900
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 14).Value = ""

$ws.Cells.Item(106, 8).Value = 27056
$ws.Cells.Item(106, 10).Value = 27056
$ws.Cells.Item(106, 12).Value = 27056
$ws.Cells.Item(106, 14).Value = -29580

$ws.Cells.Item(113, 8).Value = 770420.9
$ws.Cells.Item(113, 9).Value = 1429358.8
$ws.Cells.Item(113, 10).Value = 1660
$ws.Cells.Item(113, 11).Value = 1429358.8
$ws.Cells.Item(113, 12).Value = 1660
$ws.Cells.Item(113, 13).Value = -1427188.8
$ws.Cells.Item(113, 14).Value = -6000

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1344.5625
$ws.Cells.Item(82, 9).Value = 1075.7273
$ws.Cells.Item(82, 10).Value = 1936
$ws.Cells.Item(82, 11).Value = 1075.7273
$ws.Cells.Item(82, 12).Value = 1936
$ws.Cells.Item(82, 13).Value = -714.7273
$ws.Cells.Item(82, 14).Value = -2658

$ws.Cells.Item(85, 8).Value = 1344.5625
$ws.Cells.Item(85, 9).Value = 1075.7273
$ws.Cells.Item(85, 10).Value = 1936
$ws.Cells.Item(85, 11).Value = 1075.7273
$ws.Cells.Item(85, 12).Value = 1936
$ws.Cells.Item(85, 13).Value = 172.2727
$ws.Cells.Item(85, 14).Value = -4432

$ws.Cells.Item(93, 8).Value = 2358.5715
$ws.Cells.Item(93, 9).Value = 2500
$ws.Cells.Item(93, 10).Value = 2335
$ws.Cells.Item(93, 11).Value = 2500
$ws.Cells.Item(93, 12).Value = 2335
$ws.Cells.Item(93, 13).Value = -1252
$ws.Cells.Item(93, 14).Value = -4831

$ws.Cells.Item(100, 8).Value = 1488.3784
$ws.Cells.Item(100, 9).Value = 1395.0769
$ws.Cells.Item(100, 10).Value = 1708.909
$ws.Cells.Item(100, 11).Value = 1395.0769
$ws.Cells.Item(100, 12).Value = 1708.909
$ws.Cells.Item(100, 13).Value = -854.0769
$ws.Cells.Item(100, 14).Value = -2790.909

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 910
$ws.Cells.Item(14, 9).Value = 1000
$ws.Cells.Item(14, 10).Value = 775
$ws.Cells.Item(14, 11).Value = 1000
$ws.Cells.Item(14, 12).Value = 775
$ws.Cells.Item(14, 13).Value = -832
$ws.Cells.Item(14, 14).Value = -1111

$ws.Cells.Item(15, 8).Value = 3950
$ws.Cells.Item(15, 10).Value = 3950
$ws.Cells.Item(15, 12).Value = 3950
$ws.Cells.Item(15, 14).Value = -4526

$ws.Cells.Item(96, 8).Value = 1209.3334
$ws.Cells.Item(96, 9).Value = 820
$ws.Cells.Item(96, 10).Value = 1520.8
$ws.Cells.Item(96, 11).Value = 820
$ws.Cells.Item(96, 12).Value = 1520.8
$ws.Cells.Item(96, 13).Value = 553
$ws.Cells.Item(96, 14).Value = -4266.8
